# Update the CDA Logical model metadata for ST.r2b
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# 1. Insert a new "Jurisdiction" property row (with empty value) right after "Contact" (row 10),
#    pushing "Description" and everything below it down by one row.
#    Copy formatting from the row above so the new row matches the existing table style.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").Insert()
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# 2. Bump the Version string.
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 3. Bump the Date timestamp.
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"
